# Lume: brightness adjustment implemented.
# Rework Sheet1 from the old "register schema / PWM / hour clock" scratch
# calculations into a focused "ADC" brightness calculator.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Wipe out everything from row 8 down (old SPI/PWM block in A8:C16 and
# the old "Hour/Hyperminute" + timer scratch block in columns D..I,
# rows 20:44). The top block (rows 1-4, Crystal freq/Prescaler/etc.)
# stays untouched.
# ---------------------------------------------------------------------
$ws.Rows("8:44").Delete()

# ---------------------------------------------------------------------
# New block: timer/OVF freq (rows 8-11) - keeps the old "Divisor",
# "Timer input freq", "ICR", "OVF freq" labels, but the source cell
# that timer input freq depended on was removed, so it now errors
# with #REF! (matches the real edit exactly).
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Divisor"
$ws.Range("B8").Value = 1

$ws.Range("A9").Value = "Timer input freq"
$ws.Range("B9").Formula = "=#REF!*1000000/B8"
$ws.Range("C9").Value = "Hz"

$ws.Range("A10").Value = "ICR"
$ws.Range("B10").Value = 255

$ws.Range("A11").Value = "OVF freq"
$ws.Range("B11").Formula = "=B9/B10"
$ws.Range("C11").Value = "Hz"

# ---------------------------------------------------------------------
# New "ADC" section header (row 15), styled the same way the old
# "Register-based schema" section header was (Accent1 row style).
# ---------------------------------------------------------------------
$ws.Rows(15).Style = "Accent1"
$ws.Range("A15").Value = "ADC"

$ws.Range("A16").Value = "Reference"
$ws.Range("B16").Value = 5000
$ws.Range("C16").Value = "mV"

$ws.Range("A17").Value = "ADC MAX"
$ws.Range("B17").Value = 255
$ws.Range("C17").Value = "steps"

$ws.Range("A18").Value = "ADC step"
$ws.Range("A18").Style = "Neutral"
$ws.Range("B18").Formula = "=B16/B17"
$ws.Range("B18").Style = "Neutral"
$ws.Range("B18").NumberFormat = "0.0"
$ws.Range("C18").Value = "mv/step"

$ws.Range("A19").Value = "U"
$ws.Range("B19").Value = 3900
$ws.Range("C19").Value = "mV"

$ws.Range("A20").Value = "ADC value"
$ws.Range("A20").Style = "Neutral"
$ws.Range("B20").Formula = "=B19/B18"
$ws.Range("B20").Style = "Neutral"
$ws.Range("C20").Value = "steps"

$ws.Range("A22").Value = "Steps"
$ws.Range("B22").Value = 200
$ws.Range("C22").Value = "steps"

$ws.Range("A23").Value = "U"
$ws.Range("A23").Style = "Neutral"
$ws.Range("B23").Formula = "=B22*B18"
$ws.Range("B23").Style = "Neutral"
$ws.Range("B23").NumberFormat = "0"
$ws.Range("C23").Value = "mV"

# ---------------------------------------------------------------------
# Restore the sheet view: the old sheet scrolled down to row 12 with
# C44 selected; now it's back at the top with B28 selected (matching
# the new, much shorter sheet).
# ---------------------------------------------------------------------
$ws.Application.Goto($ws.Range("A1"), $true)
$ws.Range("B28").Select()
